# pipeflow.xlsx edit
# - Fix the h>r (Case Study 2) block: C21 was 0.99 (nonsensical, h > D = 0.6 so
#   D21=C19-C21 went negative) -> 0.5, and C22 referenced the wrong cell
#   (D21 instead of C21), which produced ACOS() of an out-of-domain argument
#   and cascaded #NUM! errors through C23:C25, rows 26-30 and H31/Q32.
#   Correcting both makes every dependent formula recalc to a real number.
# - Add a small scratch check in column L (rows 21-24) that recomputes the
#   circle area two different ways as a sanity cross-check.
# - Move/resize the second diagram picture and update the view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the buggy inputs driving the #NUM! cascade ---
$ws.Range("C21").Value = 0.5
$ws.Range("C22").Formula = "=2*ACOS((C20-C21)/C20)"

# --- New cross-check cells in column L ---
$ws.Range("L21").Formula = "=PI()*0.3^2"
$ws.Range("L22").Formula = "=2*L21"
$ws.Range("L23").Formula = "=L22/PI()"
$ws.Range("L24").Formula = "=SQRT(L23)"

# --- Reposition / resize the second picture (Picture 5) ---
$pic = $ws.Shapes.Item("Picture 5")
$pic.Left = 3508813 / 12700
$pic.Top = 2075793 / 12700
$pic.Width = 2930724 / 12700
$pic.Height = 1549414 / 12700

# --- View state: zoom + selection ---
$excel.ActiveWindow.Zoom = 145
$ws.Range("L25").Select()
